$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'63.900.65"
$ws.Range("E2").Formula = "'  +1.85%  "

$ws.Range("D3").Formula = "'2.500.58"
$ws.Range("E3").Formula = "'  +1.71%  "

$ws.Range("D4").Formula = "'0.998"
$ws.Range("E4").Formula = "'  -0.20%  "

$ws.Range("D5").Formula = "'579.31"
$ws.Range("E5").Formula = "'  +0.98%  "

$ws.Range("D6").Formula = "'151.00"
$ws.Range("E6").Formula = "'  +3.84%  "

$ws.Range("D7").Formula = "'0.998"
$ws.Range("E7").Formula = "'  -0.22%  "

$ws.Range("D8").Formula = "'0.538"
$ws.Range("E8").Formula = "'  +0.45%  "

$ws.Range("D9").Formula = "'0.112"
$ws.Range("E9").Formula = "'  +1.29%  "

$ws.Range("E10").Formula = "'  +0.25%  "

$ws.Range("D11").Formula = "'5.24"
$ws.Range("E11").Formula = "'  -0.15%  "

$ws.Range("D12").Formula = "'0.353"
$ws.Range("E12").Formula = "'  -0.53%  "

$ws.Range("D13").Formula = "'29.64"
$ws.Range("E13").Formula = "'  +3.12%  "

$ws.Range("D14").Formula = "'0.0000179"
$ws.Range("E14").Formula = "'  +1.72%  "

$ws.Range("D15").Formula = "'2.950.71"
$ws.Range("E15").Formula = "'  +1.55%  "

$ws.Range("D16").Formula = "'63.657.54"
$ws.Range("E16").Formula = "'  +1.32%  "

$ws.Range("D17").Formula = "'2.498.72"
$ws.Range("E17").Formula = "'  +1.64%  "

$ws.Range("D18").Formula = "'7.81"
$ws.Range("E18").Formula = "'  -2.21%  "

$ws.Range("D19").Formula = "'10.92"
$ws.Range("E19").Formula = "'  -0.48%  "

$ws.Range("D20").Formula = "'4.24"
$ws.Range("E20").Formula = "'  +2.60%  "

$ws.Range("D21").Formula = "'2.28"
$ws.Range("E21").Formula = "'  +2.75%  "

$ws.Range("D22").Formula = "'326.16"
$ws.Range("E22").Formula = "'  -0.11%  "

$ws.Range("E23").Formula = "'  -0.04%  "

$ws.Range("D24").Formula = "'10.16"
$ws.Range("E24").Formula = "'  +1.76%  "

$ws.Range("D25").Formula = "'671.16"
$ws.Range("E25").Formula = "'  +2.94%  "

$ws.Range("D26").Formula = "'65.22"
$ws.Range("E26").Formula = "'  -0.62%  "

$ws.Range("D27").Formula = "'0.0000100"
$ws.Range("E27").Formula = "'  +3.13%  "

$ws.Range("D28").Formula = "'2.599.78"
$ws.Range("E28").Formula = "'  +0.67%  "

$ws.Range("B29").Formula = "'Fetch.AI"
$ws.Range("C29").Formula = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").Formula = "'1.49"
$ws.Range("E29").Formula = "'  +2.77%  "

$ws.Range("B30").Formula = "'Binance-PegBSC-USD"
$ws.Range("C30").Formula = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Formula = "'0.993"
$ws.Range("E30").Formula = "'  -0.54%  "

$ws.Range("D31").Formula = "'8.01"
$ws.Range("E31").Formula = "'  +0.35%  "

$ws.Range("D32").Formula = "'1.85"
$ws.Range("E32").Formula = "'  +0.35%  "

$ws.Range("D33").Formula = "'0.135"
$ws.Range("E33").Formula = "'  +1.20%  "

$ws.Range("D34").Formula = "'0.997"
$ws.Range("E34").Formula = "'  -0.24%  "

$ws.Range("D35").Formula = "'1.54"
$ws.Range("E35").Formula = "'  +1.59%  "

$ws.Range("D36").Formula = "'4.79"
$ws.Range("E36").Formula = "'  +0.90%  "

$ws.Range("D37").Formula = "'5.55"
$ws.Range("E37").Formula = "'  +3.54%  "

$ws.Range("D38").Formula = "'0.370"
$ws.Range("E38").Formula = "'  +0.63%  "

$ws.Range("D39").Formula = "'151.70"
$ws.Range("E39").Formula = "'  -0.55%  "

$ws.Range("D40").Formula = "'18.73"
$ws.Range("E40").Formula = "'  +0.40%  "

$ws.Range("D41").Formula = "'2.80"
$ws.Range("E41").Formula = "'  +1.79%  "

$ws.Range("D42").Formula = "'1.77"
$ws.Range("E42").Formula = "'  +3.23%  "

$ws.Range("D43").Formula = "'0.998"
$ws.Range("E43").Formula = "'  -0.07%  "

$ws.Range("D44").Formula = "'159.09"
$ws.Range("E44").Formula = "'  +4.12%  "

$ws.Range("D45").Formula = "'0.0₆0300"
$ws.Range("E45").Formula = "'  -4.74%  "

$ws.Range("D46").Formula = "'15.43"
$ws.Range("E46").Formula = "'  +1.42%  "

$ws.Range("D47").Formula = "'3.62"
$ws.Range("E47").Formula = "'  +1.25%  "

$ws.Range("D48").Formula = "'20.97"
$ws.Range("E48").Formula = "'  +3.61%  "

$ws.Range("D49").Formula = "'0.614"
$ws.Range("E49").Formula = "'  +1.42%  "

$ws.Range("D50").Formula = "'0.0517"
$ws.Range("E50").Formula = "'  +1.13%  "

$ws.Range("D51").Formula = "'0.0913"
$ws.Range("E51").Formula = "'  +0.00%  "
